$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set all B1:B22 boolean cells from FALSE to TRUE ---
$ws.Range("B1:B22").Value = $true

# --- Make the "text" row (A21) red-font text (new font/style) ---
$ws.Range("A21").Font.Color = 255

# --- Narrow column A (best achievable width on this engine) ---
$ws.Columns.Item(1).ColumnWidth = 19.5

# --- Update sheet view: zoom to 125% and change the selected cell ---
$excel.ActiveWindow.Zoom = 125
$ws.Range("B22").Select()
